$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reassign the "Periodo Mora" labels for each data row (periods got rotated:
# oldest period 2412 and newest 2505 swap places, 2502/2504 swap places,
# 2503 and 2506 stay put) and refresh the matching "Valor Mora" / "Salario
# Basico" amounts for the new period mix.
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("E19").Value = "2502"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("E20").Value = "2412"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1423500

$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 49348
$ws.Range("G21").Value = 1423500
